$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.65"
$ws.Range("D3").Value = "'21.98"
$ws.Range("D4").Value = "'5.409"
$ws.Range("D5").Value = "'0.05782"
$ws.Range("D6").Value = "'3.388"
$ws.Range("D7").Value = "'6.330"
$ws.Range("D8").Value = "'0.8095"
$ws.Range("D9").Value = "'0.9563"
$ws.Range("D10").Value = "'0.1425"
$ws.Range("D11").Value = "'0.07498"
$ws.Range("D12").Value = "'0.03186"
$ws.Range("D13").Value = "'0.03014"
$ws.Range("D14").Value = "'4.129"
$ws.Range("D15").Value = "'0.09413"
$ws.Range("D16").Value = "'0.001591"
$ws.Range("D17").Value = "'0.04813"
$ws.Range("D18").Value = "'0.0005891"
$ws.Range("D19").Value = "'0.006184"
$ws.Range("D20").Value = "'0.004109"
$ws.Range("D21").Value = "'0.0009963"
$ws.Range("D22").Value = "'0.0001499"
$ws.Range("D23").Value = "'3.766"
$ws.Range("D24").Value = "'2.231"
$ws.Range("D25").Value = "'0.3229"
$ws.Range("D27").Value = "'0.0001290"
$ws.Range("D40").Value = "'0.03882"
$ws.Range("D41").Value = "'0.006325"
$ws.Range("D42").Value = "'0.1074"
$ws.Range("D43").Value = "'0.002997"
$ws.Range("D44").Value = "'0.006082"
$ws.Range("D45").Value = "'0.00005594"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.3801"
$ws.Range("D48").Value = "'0.1427"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D50").Value = "'0.01009"
